{"js": "// Remove the stray \" Target Audience\\t4\" TOC entry text that was appended\n// after the TOC field's closing fldChar (fldCharType=\"end\") in the last\n// paragraph of the table of contents. The paragraph itself (and its\n// fldChar \"end\" run) must remain untouched \u2014 only the trailing\n// space / \"Target Audience\" / tab / \"4\" runs are deleted.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the TOC paragraph that still carries the orphaned entry text.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text === \" Target Audience\\t4\") {\n    target = p;\n    break;\n  }\n}\n\nif (target) {\n  // Search within that paragraph for the exact trailing text (space +\n  // \"Target Audience\" + tab + \"4\") and delete just that range, leaving the\n  // preceding fldChar \"end\" run (and its formatting) intact.\n  const hits = target.search(\" Target Audience\\t4\", { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length > 0) {\n    hits.items[0].delete();\n    await context.sync();\n  }\n}\n", "ps1": "# Remove the stray \" Target Audience`t4\" TOC entry text that was appended\n# after the TOC field's closing fldChar (fldCharType=\"end\") in the last\n# paragraph of the table of contents. The paragraph itself (and its\n# fldChar \"end\" run) must remain untouched - only the trailing\n# space / \"Target Audience\" / tab / \"4\" text is deleted.\n\n$d = $word.ActiveDocument\n\n$search = $d.Content\n$search.Find.ClearFormatting()\n$search.Find.MatchCase = $false\n$search.Find.MatchWholeWord = $false\n$found = $search.Find.Execute(\"Target Audience\")\n\nif ($found) {\n    # The paragraph that contains the hit is the TOC paragraph carrying the\n    # orphaned entry.\n    $para = $search.Paragraphs(1)\n    $delRange = $para.Range\n\n    # Skip the paragraph's first character - that's the (empty) run holding\n    # the TOC field's closing fldChar - and delete everything else in the\n    # paragraph (the space, \"Target Audience\", the tab, and the page number),\n    # leaving the paragraph mark and the fldChar run untouched.\n    $delRange.MoveStart(1, 1)\n    $delRange.Delete()\n}\n"}
